# Add "sales page" data sheet (Category -> Brand mapping) as a new Sheet2,
# placed after the existing Sheet1, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 so tab order is Sheet1, Sheet2.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 19.109375
$ws.Columns.Item(2).ColumnWidth = 15.5546875
$ws.Columns.Item(3).ColumnWidth = 12.33203125
$ws.Columns.Item(4).ColumnWidth = 17.5546875

# ---- Data ----
$data = @(
  @("Category_id","categoryName","brandIds","brandName"),
  @("cat_001","AC","brand_001,brand_010","Godrej, Havells"),
  @("cat_002","Laptop","brand_002","Samsung"),
  @("cat_003","Others","brand_001","Godrej"),
  @("cat_004","Refrigerator","brand_001,brand_010","Godrej, Havells"),
  @("cat_005","SmartWatch","brand_002","Samsung"),
  @("cat_006","Smartphone","brand_002","Samsung"),
  @("cat_007","Tab","brand_002","Samsung"),
  @("cat_008","Washing Machine","brand_001,brand_010","Godrej, Havells")
)

for ($r = 0; $r -lt $data.Length; $r++) {
  for ($c = 0; $c -lt $data[$r].Length; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
  }
}

# ---- Row heights ----
$tallRows = @(2, 5, 9)
for ($r = 1; $r -le 9; $r++) {
  if ($tallRows -contains $r) {
    $ws.Rows.Item($r).RowHeight = 42
  } else {
    $ws.Rows.Item($r).RowHeight = 28.2
  }
}

# ---- Header row style: bold Arial, light-grey fill, medium light-grey border,
#      left/center aligned, wrapped text ----
$header = $ws.Range("A1:D1")
$header.Font.Name = "Arial"
$header.Font.Bold = $true
$header.Interior.Color = 15921906
$header.HorizontalAlignment = -4131
$header.VerticalAlignment = -4108
$header.WrapText = $true
$header.Borders.Weight = -4138
$header.Borders.Color = 14540253

# ---- Data rows style: regular Arial, medium light-grey border,
#      left/center aligned, wrapped text ----
$bodyRange = $ws.Range("A2:D9")
$bodyRange.Font.Name = "Arial"
$bodyRange.HorizontalAlignment = -4131
$bodyRange.VerticalAlignment = -4108
$bodyRange.WrapText = $true
$bodyRange.Borders.Weight = -4138
$bodyRange.Borders.Color = 14540253

# ---- Selection / active sheet ----
$ws.Range("K13").Select() | Out-Null
